# Attendance update: record actual headcounts for the first Fall 2024
# seminar (the "women in APEC" / Agriculture, Labor, Gender talk) on the
# "2024 - Fall" sheet. Column G = in-person attendees, column H = virtual
# attendees. Everything else (per-row running total in I/K, the Total /
# Average / Percent-of-capacity summary rows 17-19, and the Attendance
# Descriptives roll-up sheet) is formula-driven and recalculates
# automatically once the inputs are entered.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024 - Fall")

$ws.Range("G2").Value = 14
$ws.Range("H2").Value = 2

# Leave the selection where the user ended up after typing the figures in.
$ws.Activate()
$ws.Range("I4").Select()
